$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-06-22 Saturday" "2024-06-23 Sunday"

Replace-Text "11×90=" "16×37="
Replace-Text "67×80=" "85×59="
Replace-Text "78×85=" "52×32="
Replace-Text "51×43=" "11×11="
Replace-Text "21×58=" "52×92="
Replace-Text "93×46=" "24×65="
Replace-Text "64×57=" "35×26="
Replace-Text "92×97=" "95×36="
Replace-Text "14×15=" "99×94="
Replace-Text "64×68=" "22×73="
Replace-Text "63×52=" "30×24="
Replace-Text "33×43=" "86×68="
Replace-Text "40×25=" "33×64="
Replace-Text "85×19=" "76×41="
Replace-Text "63×11=" "45×74="
Replace-Text "30×90=" "88×39="
Replace-Text "49×33=" "56×46="
Replace-Text "34×61=" "84×82="
Replace-Text "92×65=" "42×94="
Replace-Text "78×61=" "62×93="
Replace-Text "37×41=" "46×12="
Replace-Text "56×86=" "65×72="
Replace-Text "38×13=" "71×50="
Replace-Text "45×66=" "39×73="
Replace-Text "93×23=" "81×48="
